# Update "想去人数" (F column) counts on 展览 (sheet1), 演出 (sheet2),
# and 全部类型 (sheet4) worksheets to match the refreshed scrape snapshot.
# 本地生活 (sheet3) is untouched by this update.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsShow    = $wb.Worksheets.Item("演出")
$wsAll     = $wb.Worksheets.Item("全部类型")

# 展览 (Exhibitions)
$wsExhibit.Range("F6").Value  = 208
$wsExhibit.Range("F13").Value = 96
$wsExhibit.Range("F14").Value = 2058
$wsExhibit.Range("F16").Value = 19
$wsExhibit.Range("F18").Value = 476
$wsExhibit.Range("F19").Value = 147
$wsExhibit.Range("F20").Value = 73
$wsExhibit.Range("F23").Value = 1550
$wsExhibit.Range("F24").Value = 3752
$wsExhibit.Range("F28").Value = 1121
$wsExhibit.Range("F29").Value = 105
$wsExhibit.Range("F30").Value = 1904
$wsExhibit.Range("F32").Value = 460
$wsExhibit.Range("F33").Value = 67
$wsExhibit.Range("F34").Value = 278
$wsExhibit.Range("F35").Value = 408
$wsExhibit.Range("F37").Value = 651
$wsExhibit.Range("F39").Value = 381

# 演出 (Shows)
$wsShow.Range("F2").Value = 14

# 全部类型 (All types)
$wsAll.Range("F6").Value  = 208
$wsAll.Range("F13").Value = 96
$wsAll.Range("F14").Value = 2058
$wsAll.Range("F16").Value = 14
$wsAll.Range("F17").Value = 19
$wsAll.Range("F19").Value = 476
$wsAll.Range("F20").Value = 147
$wsAll.Range("F21").Value = 73
$wsAll.Range("F24").Value = 1550
$wsAll.Range("F25").Value = 3752
$wsAll.Range("F29").Value = 1121
$wsAll.Range("F30").Value = 105
$wsAll.Range("F31").Value = 1905
$wsAll.Range("F33").Value = 460
$wsAll.Range("F34").Value = 67
$wsAll.Range("F35").Value = 278
$wsAll.Range("F36").Value = 408
$wsAll.Range("F38").Value = 651
$wsAll.Range("F40").Value = 381
